$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Experience the beach with Bondi Break slot. Play for free with engaging gameplay, potential wins of up to 6,250x your stake and a lively beach-themed design.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) Near the end of the document: drop the duplicated bold title
#    paragraph, and replace the italic meta-description paragraph's
#    text with the new feature-image prompt.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($count - 1)
$dupTitlePara.Range.Delete()

$count2 = $d.Paragraphs.Count
$imgPromptPara = $d.Paragraphs($count2)
$fullRange = $imgPromptPara.Range
$textRange = $d.Range($fullRange.Start, $fullRange.End - 1)

$newPrompt = 'Create an eye-catching feature image for &quot;Bondi Break&quot; that captures the beach theme of the game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be holding a surfboard and standing on a beach with waves in the background. The words &quot;Bondi Break&quot; should be written in bold, colorful letters above the warrior. The image should be bright and playful with a fun and adventurous vibe that captures the essence of the game.'
$promptXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>' + $newPrompt + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$textRange.InsertXML($promptXml)

Write-Host "Edit complete"
